# feat: add 2022-Q1 data
#
# Adds a new "2022-Q1" worksheet (per-fund holding breakdown) right before
# the "总计" (totals) summary sheet, and inserts a corresponding summary
# row at the top of "总计".

function Set-TextValue {
    # Write $value to $range as TEXT (keeps leading zeros / decimal-looking
    # strings like "008809" or "15.09" from being coerced to numbers),
    # without leaving a stray explicit Text number-format behind.
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) New "2022-Q1" sheet, inserted immediately before "总计"
# ---------------------------------------------------------------------
# NOTE: fetch the "总计" handle fresh and pass it straight into Add() -
# a handle captured earlier and reused *after* Add() ends up referring
# to the newly inserted sheet instead, so don't stash it in a variable
# that outlives the Add() call.
$funds = $wb.Worksheets.Add($wb.Worksheets.Item("总计"))
$funds.Name = "2022-Q1"

# Seed layout/styles from the structurally-identical "2021-Q1" sheet
# (same headers, same per-fund-row styling) so fonts/borders/alignment
# match the other quarter sheets instead of Excel's blank-sheet defaults.
$template = $wb.Worksheets.Item("2021-Q1")
$template.Range("A1:H5").Copy($funds.Range("A1"))
$funds.Range("A1").ClearContents()
# Template only has 4 data rows (A1:H5) - this quarter needs 5, so clone
# the last data row to make room for a 5th (style follows along).
$funds.Range("A5:H5").Copy($funds.Range("A6"))

# -- header row (unchanged text, just making sure it reads correctly) --
$funds.Range("B1").Value = "基金代码"
$funds.Range("C1").Value = "基金名称"
$funds.Range("D1").Value = "基金规模"
$funds.Range("E1").Value = "股票总仓位"
$funds.Range("F1").Value = "仓位占比"
$funds.Range("G1").Value = "持有市值(亿元)"
$funds.Range("H1").Value = "仓位排名"

# -- data rows --
$fundRows = @(
    @{ A=0; B="008809"; C="安信民稳增长混合A";       D="15.09"; E="44.77"; F="2.21"; G="0.3335"; H=7 },
    @{ A=1; B="012256"; C="安信丰穗一年持有混合A";     D="26.49"; E="20.42"; F="1.17"; G="0.3099"; H=7 },
    @{ A=2; B="009849"; C="安信稳健聚申一年持有期混合"; D="12.33"; E="30.81"; F="1.91"; G="0.2355"; H=6 },
    @{ A=3; B="008810"; C="安信民稳增长混合C";       D="6.76";  E="44.77"; F="2.21"; G="0.1494"; H=7 },
    @{ A=4; B="012257"; C="安信丰穗一年持有混合C";     D="2.41";  E="20.42"; F="1.17"; G="0.0282"; H=7 }
)

$r = 2
foreach ($row in $fundRows) {
    $funds.Range("A$r").Value = $row.A
    Set-TextValue $funds.Range("B$r") $row.B
    Set-TextValue $funds.Range("C$r") $row.C
    Set-TextValue $funds.Range("D$r") $row.D
    Set-TextValue $funds.Range("E$r") $row.E
    Set-TextValue $funds.Range("F$r") $row.F
    Set-TextValue $funds.Range("G$r") $row.G
    $funds.Range("H$r").Value = $row.H
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2) Prepend a "2022-Q1" row to the "总计" summary sheet
# ---------------------------------------------------------------------
# Fetch a fresh handle (see note above - this is the first use after the
# Worksheets.Add() call, so it must be re-resolved by name here).
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

# The inherited formatting on the freshly-inserted row doesn't match the
# plain (unstyled) data cells used elsewhere in this table - strip it.
$total.Range("B2:D2").ClearFormats()

$total.Range("A2").Value = 0
Set-TextValue $total.Range("B2") "2022-Q1"
$total.Range("C2").Value = 5
$total.Range("D2").Value = 1.06

# Restore the index-column style (border/bold/center) that row 2 should
# carry, matching every other row in column A.
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)  # xlPasteFormats

# Renumber the (now shifted-down) running index in column A so it still
# reads 0,1,2,3,4,5 top to bottom.
$idx = 1
for ($row = 3; $row -le 7; $row++) {
    $total.Range("A$row").Value = $idx
    $idx = $idx + 1
}
